# dynamic_programming.pptx: mark the three quiz "Smiley Face" answer shapes
# with alt text "QuizAnswer" so accessibility/automation tooling can find
# the answer shapes that get revealed on click.
$p = $ppt.ActivePresentation

# Slide 5 (sldId 765): "Smiley Face 3" -> Content Placeholder quiz slide.
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item("Smiley Face 3").AlternativeText = "QuizAnswer"

# Slide 11 (sldId 771): "Smiley Face 4" and "Smiley Face 5" -> quiz table slide.
$s11 = $p.Slides.Item(11)
$s11.Shapes.Item("Smiley Face 4").AlternativeText = "QuizAnswer"
$s11.Shapes.Item("Smiley Face 5").AlternativeText = "QuizAnswer"
